# Marks.xlsx / Sheet1 — resolve composite-key integrity issue on the
# Internals/Externals exam rows and widen the "Exam Name" column (C) to
# fit the new longer labels ("First Internals" / "Second Internals").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: first internal exam, subject 10cs42 -----------------------
$ws.Range("B2").Value = "Internals"
$ws.Range("C2").Value = "First Internals"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "a"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "10cs42"
$ws.Range("H2").Value = 2014
$ws.Range("I2").Value = "1ru353"

# --- Row 3: second internal exam, same subject 10cs42 ------------------
$ws.Range("B3").Value = "Internals"
$ws.Range("C3").Value = "Second Internals"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "a"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "10cs42"
$ws.Range("H3").Value = 2014
$ws.Range("I3").Value = "1ru353"

# --- Row 4: first internal exam, subject 10cs49 -------------------------
$ws.Range("B4").Value = "Internals"
$ws.Range("C4").Value = "First Internals"
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = "a"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "10cs49"
$ws.Range("H4").Value = 2014
$ws.Range("I4").Value = "1ru353"

# --- Row 5: externals, subject 10cs50 -----------------------------------
$ws.Range("B5").Value = "Externals"
$ws.Range("C5").Value = "Externals"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = "a"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "10cs50"
$ws.Range("H5").Value = 2014
$ws.Range("I5").Value = "1ru353"

# --- Row 6: externals, subject 10cs48 -----------------------------------
$ws.Range("B6").Value = "Externals"
$ws.Range("C6").Value = "Externals"
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = "a"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "10cs48"
$ws.Range("H6").Value = 2014
$ws.Range("I6").Value = "1ru353"

# --- Widen column C (Exam Name) so the longer labels fit ----------------
$ws.Columns.Item(3).ColumnWidth = 14.333333333333332

# --- Move the active selection to N1, matching the saved view ----------
$ws.Range("N1").Select()
